$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge the old "GPA" and "Certificate" strings into a single cell with the
# LaTeX \href markup, then remove the now-empty row that used to hold the
# certificate link.
$ws.Range("E2").Value = "GPA: 97/100 (see \href{https://www.coursera.org/account/accomplishments/verify/DC7ULMJ3CZWM}{certificate})"

# Delete entire row 3 (the row that only contained the certificate text),
# shifting row 4 (Statistical Programming... / Dundee) up into row 3.
$ws.Rows("3").Delete()

# Update the active selection to match the target workbook.
$ws.Range("C12").Select() | Out-Null
